# Auto-generated script applying the Sagittarius_Profits.xlsx market-data refresh diff
# Updates specific H:N (price/profit) cells across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2204.8235
$ws.Range("I40").Value = 2132.0833
$ws.Range("K40").Value = 2132.0833
$ws.Range("M40").Value = -1957.0833
$ws.Range("H70").Value = 15092.857
$ws.Range("I70").Value = 1266.6666
$ws.Range("J70").Value = 18863.637
$ws.Range("K70").Value = 3799.9998
$ws.Range("L70").Value = 56590.91099999999
$ws.Range("M70").Value = -3529.9998
$ws.Range("N70").Value = -57130.91099999999
$ws.Range("H73").Value = 15092.857
$ws.Range("I73").Value = 1266.6666
$ws.Range("J73").Value = 18863.637
$ws.Range("K73").Value = 3799.9998
$ws.Range("L73").Value = 56590.91099999999
$ws.Range("M73").Value = -2863.9998
$ws.Range("N73").Value = -58462.91099999999
$ws.Range("H98").Value = 977.6
$ws.Range("I98").Value = 977.6
$ws.Range("K98").Value = 977.6
$ws.Range("M98").Value = 520.4
$ws.Range("H122").Value = 977.6
$ws.Range("I122").Value = 977.6
$ws.Range("K122").Value = 2932.8
$ws.Range("M122").Value = -482.8000000000002

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3437.12
$ws.Range("I32").Value = 3437.12
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3437.12
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3150.12
$ws.Range("H97").Value = 1139.4
$ws.Range("I97").Value = 1139.4
$ws.Range("K97").Value = 1139.4
$ws.Range("M97").Value = -643.4000000000001
$ws.Range("H110").Value = 5286507
$ws.Range("I110").Value = 5286507
$ws.Range("K110").Value = 5286507
$ws.Range("M110").Value = -5284462

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4749.5
$ws.Range("I105").Value = 4733
$ws.Range("J105").Value = 4799
$ws.Range("K105").Value = 4733
$ws.Range("L105").Value = 4799
$ws.Range("M105").Value = -2986
$ws.Range("N105").Value = -8293

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2498.3572
$ws.Range("I31").Value = 1566.2222
$ws.Range("J31").Value = 4176.2
$ws.Range("K31").Value = 1566.2222
$ws.Range("L31").Value = 4176.2
$ws.Range("M31").Value = -1271.2222
$ws.Range("N31").Value = -4766.2
$ws.Range("H34").Value = 2498.3572
$ws.Range("I34").Value = 1566.2222
$ws.Range("J34").Value = 4176.2
$ws.Range("K34").Value = 1566.2222
$ws.Range("L34").Value = 4176.2
$ws.Range("M34").Value = -1364.2222
$ws.Range("N34").Value = -4580.2
$ws.Range("H94").Value = 14379
$ws.Range("I94").Value = 13201.2
$ws.Range("K94").Value = 13201.2
$ws.Range("M94").Value = -12750.2
$ws.Range("H132").Value = 2326.5
$ws.Range("I132").Value = 2393.25
$ws.Range("J132").Value = 2193
$ws.Range("K132").Value = 7179.75
$ws.Range("L132").Value = 6579
$ws.Range("M132").Value = -4649.75
$ws.Range("N132").Value = -11639

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 177.33333
$ws.Range("I17").Value = 177.33333
$ws.Range("K17").Value = 531.99999
$ws.Range("M17").Value = -362.99999
$ws.Range("H34").Value = 383
$ws.Range("J34").Value = 493
$ws.Range("L34").Value = 1479
$ws.Range("N34").Value = -1647
$ws.Range("H39").Value = 6973
$ws.Range("J39").Value = 6631
$ws.Range("L39").Value = 19893
$ws.Range("N39").Value = -20481

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 67.75
$ws.Range("J2").Value = 90
$ws.Range("L2").Value = 90
$ws.Range("N2").Value = -316
$ws.Range("H70").Value = 8345.625
$ws.Range("I70").Value = 7973.25
$ws.Range("J70").Value = 8718
$ws.Range("K70").Value = 7973.25
$ws.Range("L70").Value = 8718
$ws.Range("M70").Value = -7703.25
$ws.Range("N70").Value = -9258
$ws.Range("H73").Value = 8345.625
$ws.Range("I73").Value = 7973.25
$ws.Range("J73").Value = 8718
$ws.Range("K73").Value = 7973.25
$ws.Range("L73").Value = 8718
$ws.Range("M73").Value = -7037.25
$ws.Range("N73").Value = -10590
$ws.Range("H107").Value = 2097.8
$ws.Range("I107").Value = 1390.3334
$ws.Range("J107").Value = 3159
$ws.Range("K107").Value = 1390.3334
$ws.Range("L107").Value = 3159
$ws.Range("M107").Value = 529.6666
$ws.Range("N107").Value = -6999
$ws.Range("H122").Value = 3979.0527
$ws.Range("I122").Value = 3600.2
$ws.Range("J122").Value = 5399.75
$ws.Range("K122").Value = 10800.6
$ws.Range("L122").Value = 16199.25
$ws.Range("M122").Value = -8350.599999999999
$ws.Range("N122").Value = -21099.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1098.6875
$ws.Range("I16").Value = 1048.4166
$ws.Range("J16").Value = 1249.5
$ws.Range("K16").Value = 1048.4166
$ws.Range("L16").Value = 1249.5
$ws.Range("M16").Value = -878.4166
$ws.Range("N16").Value = -1589.5
$ws.Range("H40").Value = 4999.6665
$ws.Range("I40").Value = 3749.75
$ws.Range("K40").Value = 3749.75
$ws.Range("M40").Value = -3613.75
$ws.Range("H46").Value = 126500
$ws.Range("I46").Value = 167666.67
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 167666.67
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -167478.67
$ws.Range("N46").Value = -3376
$ws.Range("H61").Value = 4924.6665
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5404
$ws.Range("H113").Value = 4924.6665
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 5437.3125
$ws.Range("I122").Value = 5071.2144
$ws.Range("K122").Value = 15213.6432
$ws.Range("M122").Value = -12763.6432
$ws.Range("H132").Value = 2201.6365
$ws.Range("I132").Value = 1767
$ws.Range("K132").Value = 5301
$ws.Range("M132").Value = -2771
$ws.Range("H136").Value = 2509.4443
$ws.Range("I136").Value = 1730.3334
$ws.Range("J136").Value = 4067.6667
$ws.Range("K136").Value = 5191.0002
$ws.Range("L136").Value = 12203.0001
$ws.Range("M136").Value = -2641.0002
$ws.Range("N136").Value = -17303.0001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1668845.6
$ws.Range("I81").Value = 2783.3333
$ws.Range("J81").Value = 3334908
$ws.Range("K81").Value = 5566.6666
$ws.Range("L81").Value = 6669816
$ws.Range("M81").Value = -4505.6666
$ws.Range("N81").Value = -6671938
$ws.Range("H84").Value = 1668845.6
$ws.Range("I84").Value = 2783.3333
$ws.Range("J84").Value = 3334908
$ws.Range("K84").Value = 27833.333
$ws.Range("L84").Value = 33349080
$ws.Range("M84").Value = -22529.333
$ws.Range("N84").Value = -33359688
$ws.Range("H132").Value = 1526.56
$ws.Range("I132").Value = 1499.3334
$ws.Range("K132").Value = 4498.0002
$ws.Range("M132").Value = -1968.0002
